# Apply updated cryptocurrency price/volume data to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.454.09'
$ws.Range("E2").Value = '  -2.65%  '
$ws.Range("D3").Value = '2.220.40'
$ws.Range("E3").Value = '  -2.63%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '110.28'
$ws.Range("E5").Value = '  -10.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '298.67'
$ws.Range("E6").Value = '  +12.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("E7").Value = '  -1.28%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("E9").Value = '  -2.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.93'
$ws.Range("E10").Value = '  -7.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0923'
$ws.Range("E11").Value = '  -2.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.65'
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.81'
$ws.Range("E13").Value = '  -3.28%  '
$ws.Range("E14").Value = '  -2.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.959'
$ws.Range("E15").Value = '  +7.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.10'
$ws.Range("E16").Value = '  -2.32%  '
$ws.Range("D17").Value = '2.558.17'
$ws.Range("E17").Value = '  -2.44%  '
$ws.Range("D18").Value = '2.241.61'
$ws.Range("E18").Value = '  -2.04%  '
$ws.Range("D19").Value = '42.685.76'
$ws.Range("E19").Value = '  -2.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.31'
$ws.Range("E20").Value = '  +4.27%  '
$ws.Range("E21").Value = '  -3.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.80'
$ws.Range("E22").Value = '  +2.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.48'
$ws.Range("E23").Value = '  +20.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.29'
$ws.Range("E24").Value = '  -6.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '229.47'
$ws.Range("E25").Value = '  -3.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.34'
$ws.Range("E26").Value = '  -1.59%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -1.83%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.68'
$ws.Range("E28").Value = '  -1.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.90'
$ws.Range("E29").Value = '  -1.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.50'
$ws.Range("E30").Value = '  -9.82%  '
$ws.Range("E31").Value = '  -1.80%  '
$ws.Range("E32").Value = '  -3.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '174.26'
$ws.Range("E33").Value = '  +1.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.07'
$ws.Range("E34").Value = '  -3.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0884'
$ws.Range("E35").Value = '  -2.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.68'
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.90'
$ws.Range("E37").Value = '  +5.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.25'
$ws.Range("E38").Value = '  +1.17%  '
$ws.Range("E39").Value = '  -2.42%  '
$ws.Range("E40").Value = '  -3.14%  '
$ws.Range("E41").Value = '  -4.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.48'
$ws.Range("E42").Value = '  -1.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.235'
$ws.Range("E43").Value = '  -1.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.79'
$ws.Range("E44").Value = '  -6.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.81'
$ws.Range("E45").Value = '  -8.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.31'
$ws.Range("E47").Value = '  -4.04%  '
$ws.Range("E48").Value = '  -2.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.32'
$ws.Range("E49").Value = '  +4.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.68'
$ws.Range("E50").Value = '  +1.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.51'
$ws.Range("E51").Value = '  -0.80%  '
